$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1618257261410788
$ws.Range("C2").Value = 0.6182572614107884
$ws.Range("J2").Value = 0.02074688796680498
$ws.Range("P2").Value = 0.1120331950207469
$ws.Range("S2").Value = 0.08713692946058091
$ws.Range("B3").Value = 0.02564102564102564
$ws.Range("C3").Value = 0.05128205128205128
$ws.Range("J3").Value = 0.01282051282051282
$ws.Range("P3").Value = 0.7243589743589743
$ws.Range("S3").Value = 0.1858974358974359
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.8378378378378378
$ws.Range("S4").Value = 0.1351351351351351
$ws.Range("B6").Value = 0.04347826086956522
$ws.Range("D6").Value = 0.02898550724637681
$ws.Range("F6").Value = 0.05314009661835749
$ws.Range("J6").Value = 0.1594202898550725
$ws.Range("O6").Value = 0.01932367149758454
$ws.Range("Q6").Value = 0.1594202898550725
$ws.Range("R6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.4251207729468599
$ws.Range("B7").Value = 0.09036144578313253
$ws.Range("D7").Value = 0.01204819277108434
$ws.Range("E7").Value = 0.006024096385542169
$ws.Range("F7").Value = 0.07228915662650602
$ws.Range("J7").Value = 0.1265060240963855
$ws.Range("O7").Value = 0.006024096385542169
$ws.Range("Q7").Value = 0.1506024096385542
$ws.Range("R7").Value = 0.0783132530120482
$ws.Range("S7").Value = 0.4578313253012048
$ws.Range("B8").Value = 0.0945054945054945
$ws.Range("D8").Value = 0.01098901098901099
$ws.Range("F8").Value = 0.04835164835164835
$ws.Range("J8").Value = 0.1054945054945055
$ws.Range("O8").Value = 0.01538461538461539
$ws.Range("Q8").Value = 0.1362637362637363
$ws.Range("R8").Value = 0.1164835164835165
$ws.Range("S8").Value = 0.4725274725274725
$ws.Range("B9").Value = 0.07392996108949416
$ws.Range("D9").Value = 0.01945525291828794
$ws.Range("F9").Value = 0.06614785992217899
$ws.Range("J9").Value = 0.1050583657587549
$ws.Range("O9").Value = 0.01556420233463035
$ws.Range("Q9").Value = 0.1361867704280156
$ws.Range("R9").Value = 0.1361867704280156
$ws.Range("S9").Value = 0.4474708171206226
$ws.Range("B10").Value = 0.08152958152958154
$ws.Range("D10").Value = 0.0165945165945166
$ws.Range("E10").Value = 0.0007215007215007215
$ws.Range("F10").Value = 0.06349206349206349
$ws.Range("J10").Value = 0.09956709956709957
$ws.Range("O10").Value = 0.01298701298701299
$ws.Range("Q10").Value = 0.2005772005772006
$ws.Range("R10").Value = 0.119047619047619
$ws.Range("S10").Value = 0.4054834054834055
$ws.Range("G11").Value = 0.13671875
$ws.Range("J11").Value = 0.09765625
$ws.Range("K11").Value = 0.1953125
$ws.Range("L11").Value = 0.546875
$ws.Range("S11").Value = 0.0234375
$ws.Range("G12").Value = 0.7210884353741497
$ws.Range("J12").Value = 0.1972789115646258
$ws.Range("K12").Value = 0.01360544217687075
$ws.Range("L12").Value = 0.006802721088435374
$ws.Range("S12").Value = 0.06122448979591837
$ws.Range("G13").Value = 0.7297297297297297
$ws.Range("J13").Value = 0.2432432432432433
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.01260504201680672
$ws.Range("H15").Value = 0.1470588235294118
$ws.Range("I15").Value = 0.1008403361344538
$ws.Range("J15").Value = 0.3823529411764706
$ws.Range("K15").Value = 0.0546218487394958
$ws.Range("M15").Value = 0.004201680672268907
$ws.Range("O15").Value = 0.05042016806722689
$ws.Range("S15").Value = 0.2478991596638656
$ws.Range("F16").Value = 0.02380952380952381
$ws.Range("H16").Value = 0.1011904761904762
$ws.Range("I16").Value = 0.09523809523809523
$ws.Range("J16").Value = 0.4523809523809524
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.0119047619047619
$ws.Range("N16").Value = 0.005952380952380952
$ws.Range("O16").Value = 0.04166666666666666
$ws.Range("S16").Value = 0.1845238095238095
$ws.Range("F17").Value = 0.006976744186046512
$ws.Range("H17").Value = 0.2
$ws.Range("I17").Value = 0.1116279069767442
$ws.Range("J17").Value = 0.4651162790697674
$ws.Range("K17").Value = 0.06511627906976744
$ws.Range("M17").Value = 0.004651162790697674
$ws.Range("O17").Value = 0.04883720930232558
$ws.Range("S17").Value = 0.09767441860465116
$ws.Range("F18").Value = 0.0103448275862069
$ws.Range("H18").Value = 0.1344827586206896
$ws.Range("I18").Value = 0.1068965517241379
$ws.Range("J18").Value = 0.4448275862068966
$ws.Range("K18").Value = 0.0896551724137931
$ws.Range("M18").Value = 0.0103448275862069
$ws.Range("O18").Value = 0.07241379310344828
$ws.Range("S18").Value = 0.1310344827586207
$ws.Range("F19").Value = 0.0102880658436214
$ws.Range("H19").Value = 0.1934156378600823
$ws.Range("I19").Value = 0.09465020576131687
$ws.Range("J19").Value = 0.3895747599451303
$ws.Range("K19").Value = 0.08230452674897119
$ws.Range("M19").Value = 0.02126200274348422
$ws.Range("N19").Value = 0.001371742112482853
$ws.Range("O19").Value = 0.08161865569272976
$ws.Range("S19").Value = 0.1255144032921811
